# Fixed climate agg date sql: column E on Sheet1 held numeric Excel date
# serials (and a couple of raw text values). Replace every Date cell (E2:E57)
# with the literal SQL expression text "to_date('...', 'yyyy-mm-dd')" and reset
# the cell back to the default (unstyled) look, since the cell is text now,
# not a formatted date.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateSqlByRow = @(
    @(2, 'to_date(''2021-11-02 '', ''yyyy-mm-dd'')'),
    @(3, 'to_date(''2022-11-06 '', ''yyyy-mm-dd'')'),
    @(4, 'to_date(''2022-06-16 '', ''yyyy-mm-dd'')'),
    @(5, 'to_date(''2021-10-29 '', ''yyyy-mm-dd'')'),
    @(6, 'to_date(''2021-06-01 '', ''yyyy-mm-dd'')'),
    @(7, 'to_date(''NULL'', ''yyyy-mm-dd'')'),
    @(8, 'to_date(''2023-11-03 '', ''yyyy-mm-dd'')'),
    @(9, 'to_date(''NULL'', ''yyyy-mm-dd'')'),
    @(10, 'to_date(''2021-07-12 '', ''yyyy-mm-dd'')'),
    @(11, 'to_date(''NULL'', ''yyyy-mm-dd'')'),
    @(12, 'to_date(''2020-04-01 '', ''yyyy-mm-dd'')'),
    @(13, 'to_date(''2020-04-01 '', ''yyyy-mm-dd'')'),
    @(14, 'to_date(''2021-10-28 '', ''yyyy-mm-dd'')'),
    @(15, 'to_date(''2021-10-28 '', ''yyyy-mm-dd'')'),
    @(16, 'to_date(''2020-12-30 '', ''yyyy-mm-dd'')'),
    @(17, 'to_date(''2021-11-21 '', ''yyyy-mm-dd'')'),
    @(18, 'to_date(''2020-12-01 '', ''yyyy-mm-dd'')'),
    @(19, 'to_date(''2019-12-01 '', ''yyyy-mm-dd'')'),
    @(20, 'to_date(''2021-07-23 '', ''yyyy-mm-dd'')'),
    @(21, 'to_date(''2021-07-23 '', ''yyyy-mm-dd'')'),
    @(22, 'to_date(''2021-08-18 '', ''yyyy-mm-dd'')'),
    @(23, 'to_date(''2021-08-18 '', ''yyyy-mm-dd'')'),
    @(24, 'to_date(''2022-08-26 '', ''yyyy-mm-dd'')'),
    @(25, 'to_date(''2022-08-26 '', ''yyyy-mm-dd'')'),
    @(26, 'to_date(''2022-09-23 '', ''yyyy-mm-dd'')'),
    @(27, 'to_date(''2022-09-23 '', ''yyyy-mm-dd'')'),
    @(28, 'to_date(''2021-06-22 '', ''yyyy-mm-dd'')'),
    @(29, 'to_date(''22/10/2021'', ''yyyy-mm-dd'')'),
    @(30, 'to_date(''2023-06-27 '', ''yyyy-mm-dd'')'),
    @(31, 'to_date(''2023-06-27 '', ''yyyy-mm-dd'')'),
    @(32, 'to_date(''2020-12-28 '', ''yyyy-mm-dd'')'),
    @(33, 'to_date(''2020-12-28 '', ''yyyy-mm-dd'')'),
    @(34, 'to_date(''2022-11-17 '', ''yyyy-mm-dd'')'),
    @(35, 'to_date(''2022-11-17 '', ''yyyy-mm-dd'')'),
    @(36, 'to_date(''2021-06-22 '', ''yyyy-mm-dd'')'),
    @(37, 'to_date(''2021-06-22 '', ''yyyy-mm-dd'')'),
    @(38, 'to_date(''2020-12-01 '', ''yyyy-mm-dd'')'),
    @(39, 'to_date(''2020-12-01 '', ''yyyy-mm-dd'')'),
    @(40, 'to_date(''2021-07-30 '', ''yyyy-mm-dd'')'),
    @(41, 'to_date(''2021-07-30 '', ''yyyy-mm-dd'')'),
    @(42, 'to_date(''2022-11-03 '', ''yyyy-mm-dd'')'),
    @(43, 'to_date(''2020-11-25 '', ''yyyy-mm-dd'')'),
    @(44, 'to_date(''2020-12-18 '', ''yyyy-mm-dd'')'),
    @(45, 'to_date(''2020-12-18 '', ''yyyy-mm-dd'')'),
    @(46, 'to_date(''2021-04-15 '', ''yyyy-mm-dd'')'),
    @(47, 'to_date(''2021-04-15 '', ''yyyy-mm-dd'')'),
    @(48, 'to_date(''2022-11-04 '', ''yyyy-mm-dd'')'),
    @(49, 'to_date(''2022-11-03 '', ''yyyy-mm-dd'')'),
    @(50, 'to_date(''2021-12-17 '', ''yyyy-mm-dd'')'),
    @(51, 'to_date(''2021-01-28 '', ''yyyy-mm-dd'')'),
    @(52, 'to_date(''2022-11-02 '', ''yyyy-mm-dd'')'),
    @(53, 'to_date(''2022-11-02 '', ''yyyy-mm-dd'')'),
    @(54, 'to_date(''2023-04-13 '', ''yyyy-mm-dd'')'),
    @(55, 'to_date(''2021-09-01 '', ''yyyy-mm-dd'')'),
    @(56, 'to_date(''2022-01-11 '', ''yyyy-mm-dd'')'),
    @(57, 'to_date(''2022-11-01 '', ''yyyy-mm-dd'')')
)

foreach ($pair in $dateSqlByRow) {
    $row = $pair[0]
    $sql = $pair[1]
    $cell = $ws.Cells.Item($row, 5)   # column E
    $cell.Value = $sql
    $cell.Style = "Normal"            # drop the yyyy-mm-dd hh:mm:ss date style
}

